# Update LR-pair stats with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1526236666666667
$ws.Range("H2").Value = 0.457871
$ws.Range("M2").Value = 27.85292233333333
$ws.Range("N2").Value = 83.55876699999999
$ws.Range("O2").Value = 0.1175699887262562
$ws.Range("P2").Value = 0.1175699887262562
$ws.Range("Q2").Value = 4.251015133895221
$ws.Range("R2").Value = 38.25913620505699
$ws.Range("S2").Value = 0.1175699887262562
$ws.Range("T2").Value = 0.1175699887262562

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1526236666666667
$ws.Range("H3").Value = 0.457871
$ws.Range("O3").Value = 0.1125380329259528
$ws.Range("P3").Value = 0.1125380329259528
$ws.Range("Q3").Value = 4.069073122231111
$ws.Range("R3").Value = 36.62165810008
$ws.Range("S3").Value = 0.1125380329259528
$ws.Range("T3").Value = 0.1125380329259528

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1526236666666667
$ws.Range("H4").Value = 0.457871
$ws.Range("M4").Value = 105.665011
$ws.Range("N4").Value = 316.995033
$ws.Range("O4").Value = 0.4460226472237104
$ws.Range("P4").Value = 0.4460226472237104
$ws.Range("Q4").Value = 16.12698141719367
$ws.Range("R4").Value = 145.142832754743
$ws.Range("S4").Value = 0.4460226472237104
$ws.Range("T4").Value = 0.4460226472237104

# Row 5 (Target cluster: MuSCs)
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1526236666666667
$ws.Range("H5").Value = 0.457871
$ws.Range("M5").Value = 4.883238666666666
$ws.Range("N5").Value = 14.649716
$ws.Range("O5").Value = 0.02061264193813266
$ws.Range("P5").Value = 0.02061264193813266
$ws.Range("Q5").Value = 0.745297790515111
$ws.Range("R5").Value = 6.707680114635999
$ws.Range("S5").Value = 0.02061264193813266
$ws.Range("T5").Value = 0.02061264193813266

# Row 6 (Target cluster: Resolving-Mac)
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1526236666666667
$ws.Range("H6").Value = 0.457871
$ws.Range("M6").Value = 71.84303666666666
$ws.Range("N6").Value = 215.52911
$ws.Range("O6").Value = 0.3032566891859479
$ws.Range("P6").Value = 0.3032566891859479
$ws.Range("Q6").Value = 10.96494768053444
$ws.Range("R6").Value = 98.68452912481
$ws.Range("S6").Value = 0.3032566891859479
$ws.Range("T6").Value = 0.3032566891859479
